$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '40.038.05'
$ws.Range("E2").Value = '  -3.99%  '
$ws.Range("D3").Value = '2.331.49'
$ws.Range("E3").Value = '  -5.71%  '
$ws.Range("E4").Value = '  -0.08%  '
Set-TextValue "D5" '307.28'
$ws.Range("E5").Value = '  -4.06%  '
Set-TextValue "D6" '85.09'
$ws.Range("E6").Value = '  -7.64%  '
Set-TextValue "D7" '0.529'
$ws.Range("E7").Value = '  -3.75%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -5.08%  '
Set-TextValue "D10" '0.0816'
$ws.Range("E10").Value = '  -4.07%  '
Set-TextValue "D11" '30.09'
$ws.Range("E11").Value = '  -8.67%  '
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = '2.690.12'
$ws.Range("E13").Value = '  -5.72%  '
$ws.Range("E14").Value = '  -6.67%  '
Set-TextValue "D15" '14.70'
$ws.Range("E15").Value = '  -4.83%  '
$ws.Range("D16").Value = '2.327.94'
$ws.Range("E16").Value = '  -5.84%  '
Set-TextValue "D17" '0.752'
$ws.Range("E17").Value = '  -4.93%  '
$ws.Range("D18").Value = '39.990.41'
$ws.Range("E18").Value = '  -3.94%  '
$ws.Range("E19").Value = '  -3.70%  '
Set-TextValue "D20" '6.09'
$ws.Range("E20").Value = '  -5.48%  '
Set-TextValue "D21" '67.59'
$ws.Range("E21").Value = '  -5.33%  '
$ws.Range("E22").Value = '  -4.84%  '
Set-TextValue "D23" '235.57'
$ws.Range("E23").Value = '  -1.63%  '
Set-TextValue "D24" '2.55'
$ws.Range("E24").Value = '  -7.30%  '
$ws.Range("E26").Value = '  -7.23%  '
Set-TextValue "D27" '23.39'
$ws.Range("E27").Value = '  -6.08%  '
Set-TextValue "D28" '2.15'
$ws.Range("E28").Value = '  -3.79%  '
$ws.Range("E29").Value = '  -5.19%  '
Set-TextValue "D30" '35.01'
$ws.Range("E30").Value = '  -4.08%  '
Set-TextValue "D31" '152.52'
$ws.Range("E31").Value = '  -2.90%  '
$ws.Range("E32").Value = '  -0.06%  '
Set-TextValue "D33" '5.12'
$ws.Range("E33").Value = '  -5.75%  '
$ws.Range("E34").Value = '  -4.53%  '
Set-TextValue "D35" '0.0723'
$ws.Range("E35").Value = '  -5.44%  '
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D38" '15.75'
$ws.Range("E38").Value = '  -7.96%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D39" '2.74'
$ws.Range("E39").Value = '  -5.02%  '
$ws.Range("E40").Value = '  -7.05%  '
Set-TextValue "D41" '3.80'
$ws.Range("E41").Value = '  -4.99%  '
Set-TextValue "D42" '2.27'
$ws.Range("E42").Value = '  -6.21%  '
$ws.Range("D43").Value = '1.936.51'
$ws.Range("E43").Value = '  -3.36%  '
$ws.Range("E44").Value = '  -5.45%  '
Set-TextValue "D45" '17.54'
$ws.Range("E45").Value = '  -5.96%  '
$ws.Range("E46").Value = '  -1.98%  '
$ws.Range("E47").Value = '  -9.44%  '
$ws.Range("D48").Value = '2.559.96'
$ws.Range("E48").Value = '  -6.21%  '
Set-TextValue "D49" '92.89'
$ws.Range("E49").Value = '  -4.67%  '
Set-TextValue "D50" '71.42'
$ws.Range("E50").Value = '  -5.56%  '
Set-TextValue "D51" '50.49'
$ws.Range("E51").Value = '  -3.03%  '
